$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$ws.Range("Z2:Z8").Value = "2025-11-13T06:53:14.541087"
$ws.Range("Z9:Z19").Value = "2025-11-13T06:53:14.542089"
$ws.Range("Z20:Z30").Value = "2025-11-13T06:53:14.543088"
$ws.Range("Z31:Z40").Value = "2025-11-13T06:53:14.544088"
$ws.Range("Z41:Z45").Value = "2025-11-13T06:53:14.545088"
$ws.Range("Z46:Z49").Value = "2025-11-13T06:53:14.756678"
$ws.Range("Z50:Z54").Value = "2025-11-13T06:53:14.757679"
$ws.Range("Z55:Z58").Value = "2025-11-13T06:53:14.758681"
$ws.Range("Z59:Z61").Value = "2025-11-13T06:53:14.759676"
$ws.Range("Z62").Value = "2025-11-13T06:53:14.760678"
$ws.Range("Z63").Value = "2025-11-13T06:53:14.762679"
$ws.Range("Z64:Z66").Value = "2025-11-13T06:53:14.765854"
$ws.Range("Z67:Z69").Value = "2025-11-13T06:53:14.766854"
$ws.Range("Z70:Z72").Value = "2025-11-13T06:53:14.767856"
$ws.Range("Z73:Z74").Value = "2025-11-13T06:53:14.768852"
$ws.Range("Z75:Z82").Value = "2025-11-13T06:53:14.983858"
$ws.Range("Z83:Z93").Value = "2025-11-13T06:53:14.984858"
$ws.Range("Z94:Z102").Value = "2025-11-13T06:53:14.985857"
